$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- Existing "직사각형 3" rectangle: move it to the left ---
$rect3 = $s.Shapes.Item(1)
$rect3.Left = 364.6530708661417   # 4631094 EMU
$rect3.Top  = 211.59173228346455  # 2687215 EMU

# --- New "직사각형 4" rectangle ---
# Duplicate the original shape twice (and discard) purely so the slide's
# internal shape-id counter advances to match the id (5) used by the
# target file, then keep the third duplicate as the real new shape.
$tmp1 = $rect3.Duplicate()
$tmp1.Item(1).Delete()

$tmp2 = $rect3.Duplicate()
$tmp2.Item(1).Delete()

$dup = $rect3.Duplicate()
$rect4 = $dup.Item(1)

$rect4.Name = "직사각형 4"
$rect4.Left = 480.0                # 6096000 EMU
$rect4.Top  = 211.59181102362206   # 2687216 EMU
$rect4.Fill.ForeColor.RGB = 0xDEDEF2   # srgbClr F2DEDE (COM RGB is 0xBBGGRR)
